# Update the "cards" column (C) values for rows 4-6, rotating them:
#   C4: "rest of the cards"        -> "manage turn order"
#   C5: "manage turn order"        -> "end game and count points"
#   C6: "end game and count points"-> "rest of the cards"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "manage turn order"
$ws.Range("C5").Value = "end game and count points"
$ws.Range("C6").Value = "rest of the cards"

# Update the active selection on the sheet (was D17, now C4).
$ws.Range("C4").Select()

# Update the workbook window X position (was 3636, now 4812).
$excel.ActiveWindow.Left = 4812
